# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.1496068669990043;  E = 0.5333859586016987 }
    3  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3993.344853322108;   E = 13.86384647080068 }
    4  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;   E = 0.5333859586016987 }
    5  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265;  E = 0.5333859586016987 }
    6  = @{ B = 0.6545652718822623; C = 1.626987699542094; D = 0.7210945179870265;  E = 13.86384647080068 }
    7  = @{ B = 0.01253208636536152;C = 0.3048912486333797;D = 3.223369029078222;   E = 0.5333859586016987 }
    8  = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;   E = 0.5333859586016987 }
    9  = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 0.1496068669990043;  E = 0.5333859586016987 }
    10 = @{ B = 0.01253208636536152;C = 0.3048912486333797;D = 0.1496068669990043;  E = 0.5333859586016987 }
    11 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 18.71679738969934;   E = 13.86384647080068 }
    12 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;   E = 0.5333859586016987 }
}

foreach ($r in $values.Keys) {
    $row = $values[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.B + $row.C + $row.D + $row.E
}
